# Insert a new weekly data row at row 146 (pushing the existing rows 146-221
# down to 147-222), then populate the new row with the latest "Acelga"
# price-report record for Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 146..221 down to 147..222, leaving a blank row 146 (inherits the
# Date-column style from the row below, same as a native Excel row insert).
$ws.Rows.Item(146).Insert()

# Populate the newly inserted row 146 with the new record.
$ws.Range("A146").Value = 10
$ws.Range("B146").Value = "Vega Modelo de Temuco"
$ws.Range("C146").Value = "La Araucanía"
$ws.Range("D146").Value = 44518
$ws.Range("E146").Value = 9
$ws.Range("F146").Value = 100112009
$ws.Range("G146").Value = "Acelga"
$ws.Range("H146").Value = "Sin especificar"
$ws.Range("I146").Value = "Primera"
$ws.Range("J146").Value = 90
$ws.Range("K146").Value = 8000
$ws.Range("L146").Value = 9000
$ws.Range("M146").Value = 8444
$ws.Range("N146").Value = "$/docena de atados (12 kilos)"
$ws.Range("O146").Value = "Provincia de Cautín"
$ws.Range("P146").Value = 704
$ws.Range("Q146").Value = 12
$ws.Range("R146").Value = "Hortaliza"
